$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 72-73; everything currently at row 72 downward
# shifts down by two rows (old 72 -> 74, ..., old 165 -> 167).
$ws.Rows("72:73").Insert()

# Populate the two newly inserted rows with the new price records.
$ws.Range("A72").Value = 7
$ws.Range("B72").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C72").Value = "Ñuble"
$ws.Range("D72").Value = 44494
$ws.Range("E72").Value = 16
$ws.Range("F72").Value = 100112002
$ws.Range("G72").Value = "Pimiento"
$ws.Range("H72").Value = "Zafiro rojo"
$ws.Range("I72").Value = "Primera"
$ws.Range("J72").Value = 100
$ws.Range("K72").Value = 43000
$ws.Range("L72").Value = 44000
$ws.Range("M72").Value = 43500
$ws.Range("N72").Value = "`$/caja 15 kilos"
$ws.Range("O72").Value = "Región de Arica y Parinacota"
$ws.Range("P72").Value = 2900
$ws.Range("Q72").Value = 15
$ws.Range("R72").Value = "Hortaliza"

$ws.Range("A73").Value = 7
$ws.Range("B73").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C73").Value = "Ñuble"
$ws.Range("D73").Value = 44494
$ws.Range("E73").Value = 16
$ws.Range("F73").Value = 100112002
$ws.Range("G73").Value = "Pimiento"
$ws.Range("H73").Value = "Zafiro verde"
$ws.Range("I73").Value = "Primera"
$ws.Range("J73").Value = 100
$ws.Range("K73").Value = 41000
$ws.Range("L73").Value = 42000
$ws.Range("M73").Value = 41500
$ws.Range("N73").Value = "`$/caja 15 kilos"
$ws.Range("O73").Value = "Región de Arica y Parinacota"
$ws.Range("P73").Value = 2767
$ws.Range("Q73").Value = 15
$ws.Range("R73").Value = "Hortaliza"
